# Fruta / hortaliza, semanal
# New weekly data (3 rows) is prepended right before the old row 27 data,
# pushing the existing rows 27-36 down to become rows 30-39.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows above the current row 27, shifting rows 27:36 down to 30:39.
$ws.Rows("27:29").Insert()

# Fill in the 3 new rows (27-29) with this week's Espárragos data.
$newRows = @(
    @{ Row=27; D=44474; H="Sin especificar"; I="Banquete"; J=180; K=1600; L=1600; M=1600; N="$/kilo";  O="Región Metropolitana"; P=1600; Q=1 },
    @{ Row=28; D=44474; H="Sin especificar"; I="Primera";  J=150; K=1400; L=1400; M=1400; N="$/kilo";  O="Región Metropolitana"; P=1400; Q=1 },
    @{ Row=29; D=44474; H="Sin especificar"; I="Segunda";  J=130; K=1200; L=1200; M=1200; N="$/kilo";  O="Región Metropolitana"; P=1200; Q=1 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = 12
    $ws.Range("B$row").Value = "Mapocho Venta Directa de Santiago"
    $ws.Range("C$row").Value = "Metropolitana"
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = 13
    $ws.Range("F$row").Value = 300000000
    $ws.Range("G$row").Value = "Espárragos"
    $ws.Range("H$row").Value = $r.H
    $ws.Range("I$row").Value = $r.I
    $ws.Range("J$row").Value = $r.J
    $ws.Range("K$row").Value = $r.K
    $ws.Range("L$row").Value = $r.L
    $ws.Range("M$row").Value = $r.M
    $ws.Range("N$row").Value = $r.N
    $ws.Range("O$row").Value = $r.O
    $ws.Range("P$row").Value = $r.P
    $ws.Range("Q$row").Value = $r.Q
    $ws.Range("R$row").Value = "Hortaliza"
}
